$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D updates are stored as literal text (prices contain
# trailing zeros / thousands separators that must not be reinterpreted
# as numbers), matching the original inline-string cell content.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "88.584.92"
$ws.Range("E2").Value = "  +8.75%  "
$ws.Range("D3").Value = "3.403.83"
$ws.Range("E3").Value = "  +7.52%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "220.07"
$ws.Range("E5").Value = "  +4.56%  "
$ws.Range("D6").Value = "649.00"
$ws.Range("E6").Value = "  +4.38%  "
$ws.Range("D7").Value = "0.424"
$ws.Range("E7").Value = "  +51.76%  "
$ws.Range("D8").Value = "0.672"
$ws.Range("E8").Value = "  +15.36%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "3.398.62"
$ws.Range("E10").Value = "  +7.61%  "
$ws.Range("D11").Value = "0.623"
$ws.Range("E11").Value = "  +7.09%  "
$ws.Range("D12").Value = "0.0000288"
$ws.Range("E12").Value = "  +15.78%  "
$ws.Range("D13").Value = "36.88"
$ws.Range("E13").Value = "  +17.32%  "
$ws.Range("D14").Value = "0.170"
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.027.30"
$ws.Range("E15").Value = "  +7.80%  "
$ws.Range("B16").Value = "Toncoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D16").Value = "5.61"
$ws.Range("E16").Value = "  +5.87%  "
$ws.Range("D17").Value = "88.439.32"
$ws.Range("E17").Value = "  +8.82%  "
$ws.Range("D18").Value = "3.408.19"
$ws.Range("E18").Value = "  +8.05%  "
$ws.Range("D19").Value = "15.10"
$ws.Range("E19").Value = "  +8.62%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "9.72"
$ws.Range("E20").Value = "  +8.51%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "461.68"
$ws.Range("E21").Value = "  +7.03%  "
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D22").Value = "3.07"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").Value = "5.69"
$ws.Range("E23").Value = "  +12.12%  "
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "7.45"
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "5.60"
$ws.Range("E25").Value = "  +6.52%  "
$ws.Range("D26").Value = "12.86"
$ws.Range("E26").Value = "  +19.39%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "3.583.72"
$ws.Range("E27").Value = "  +8.02%  "
$ws.Range("B28").Value = "Litecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D28").Value = "81.46"
$ws.Range("E28").Value = "  +6.88%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0000144"
$ws.Range("E29").Value = "  +19.86%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "0.189"
$ws.Range("E31").Value = "  +37.14%  "
$ws.Range("D32").Value = "9.53"
$ws.Range("E32").Value = "  +7.05%  "
$ws.Range("D33").Value = "584.44"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  +3.75%  "
$ws.Range("D36").Value = "2.11"
$ws.Range("E36").Value = "  +6.32%  "
$ws.Range("D37").Value = "7.51"
$ws.Range("E37").Value = "  +24.92%  "
$ws.Range("D38").Value = "0.144"
$ws.Range("E38").Value = "  -6.00%  "
$ws.Range("D39").Value = "24.10"
$ws.Range("E39").Value = "  +6.24%  "
$ws.Range("D40").Value = "0.440"
$ws.Range("E40").Value = "  +8.09%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "21.85"
$ws.Range("E41").Value = "  +5.28%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "2.13"
$ws.Range("E43").Value = "  +3.67%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "3.14"
$ws.Range("E44").Value = "  +4.69%  "
$ws.Range("D45").Value = "158.55"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "190.46"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("B48").Value = "ImmutableX"
$ws.Range("C48").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D48").Value = "1.43"
$ws.Range("E48").Value = "  +7.74%  "
$ws.Range("D49").Value = "46.74"
$ws.Range("E49").Value = "  +4.09%  "
$ws.Range("D50").Value = "4.55"
$ws.Range("E50").Value = "  +8.99%  "
$ws.Range("D51").Value = "0.676"
$ws.Range("E51").Value = "  +8.00%  "

# Restore original (default) cell formatting now that the text values
# are committed, so no stray number-format styles leak into the sheet.
$ws.Range("D2:D51").ClearFormats()
